$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update QTY column (A) values per updated test fixture data
$ws.Range("A2").Value = 23
$ws.Range("A5").Value = 0
$ws.Range("A8").Value = 0
$ws.Range("A10").Value = 0
$ws.Range("A15").Value = 0

# Move the active selection to A2
[void]$ws.Range("A2").Select()
